# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the aggregated "全部类型" sheet, reflecting a refreshed scrape.
#
# 展览 (sheet 1):
#   F3 : 2739 -> 2738
#   F5 : 19749 -> 19752
#   F7 : 2269 -> 2273
#   F8 : 750  -> 751
#   F11: 694  -> 695
#
# 全部类型 (sheet 4):
#   F8 : 2739 -> 2738
#   F10: 19749 -> 19752
#   F16: 2269 -> 2273
#   F17: 750  -> 751
#   F21: 694  -> 695

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 2738
$wsExhibition.Range("F5").Value = 19752
$wsExhibition.Range("F7").Value = 2273
$wsExhibition.Range("F8").Value = 751
$wsExhibition.Range("F11").Value = 695

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 2738
$wsAll.Range("F10").Value = 19752
$wsAll.Range("F16").Value = 2273
$wsAll.Range("F17").Value = 751
$wsAll.Range("F21").Value = 695
